$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 382
$ws1.Range("F5").Value = 415
$ws1.Range("F6").Value = 940
$ws1.Range("F9").Value = 303
$ws1.Range("F10").Value = 1222
$ws1.Range("F12").Value = 278
$ws1.Range("F13").Value = 1079
$ws1.Range("F14").Value = 423
$ws1.Range("F15").Value = 6813
$ws1.Range("F16").Value = 72
$ws1.Range("F19").Value = 7729
$ws1.Range("F20").Value = 43
$ws1.Range("F22").Value = 4164
$ws1.Range("F23").Value = 39
$ws1.Range("F24").Value = 2234
$ws1.Range("F25").Value = 952
$ws1.Range("F27").Value = 227
$ws1.Range("F31").Value = 259
$ws1.Range("F33").Value = 4
$ws1.Range("F34").Value = 1900
$ws1.Range("F36").Value = 218
$ws1.Range("F38").Value = 520
$ws1.Range("F40").Value = 1322
$ws1.Range("F41").Value = 14
$ws1.Range("F42").Value = 2016
$ws1.Range("F43").Value = 2166
$ws1.Range("F44").Value = 13

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 77
$ws2.Range("F3").Value = 61

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 247
$ws3.Range("F3").Value = 1249
$ws3.Range("F4").Value = 81

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 247
$ws4.Range("F4").Value = 1249
$ws4.Range("F5").Value = 81
$ws4.Range("F7").Value = 382
$ws4.Range("F8").Value = 415
$ws4.Range("F9").Value = 940
$ws4.Range("F11").Value = 303
$ws4.Range("F12").Value = 1222
$ws4.Range("F13").Value = 77
$ws4.Range("F14").Value = 278
$ws4.Range("F15").Value = 1079
$ws4.Range("F16").Value = 423
$ws4.Range("F17").Value = 6813
$ws4.Range("F18").Value = 72
$ws4.Range("F21").Value = 7729
$ws4.Range("F22").Value = 43
$ws4.Range("F24").Value = 4164
$ws4.Range("F25").Value = 39
$ws4.Range("F26").Value = 2234
$ws4.Range("F27").Value = 952
$ws4.Range("F29").Value = 227
$ws4.Range("F32").Value = 61
$ws4.Range("F35").Value = 259
$ws4.Range("F36").Value = 1900
$ws4.Range("F38").Value = 218
$ws4.Range("F40").Value = 520
$ws4.Range("F43").Value = 1322
$ws4.Range("F44").Value = 14
$ws4.Range("F45").Value = 2016
$ws4.Range("F47").Value = 2166
$ws4.Range("F48").Value = 13
